# Update status text and datestamps to reflect a new (out-of-sync) handback report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: not in sync with en-US"

# Status column on each sheet
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Datetime columns
$overview.Range("G2").Value = "2016-07-08 09:41:26"
$dede.Range("G2").Value = "2016-07-08 09:41:26"
$zhcn.Range("G2").Value = "2016-07-08 09:41:17"

# Widen the Status columns to fit the longer text (auto-fit result of the
# report generator after the status string grew longer). The host quantizes
# ColumnWidth to 1/6-character steps, so feed it the character width whose
# quantized result lands nearest the recorded OOXML width (33.4602203369141).
$newWidth = 32.6666666666667
$overview.Range("E1").ColumnWidth = $newWidth
$overview.Range("F1").ColumnWidth = $newWidth
$zhcn.Range("C1").ColumnWidth = $newWidth
$dede.Range("C1").ColumnWidth = $newWidth
